$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 1.845572671350878)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    4 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    6 = @(0.1190320826869504, 0.00006240767534437808, 0.1494219747398047, 0.4942365360607697, 0, 0.7627530011628691)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    8 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
